# Auto-generated Excel COM-interop script
# Updates FFXIV crafting-leve profit sheets (H:N columns) with refreshed
# market-board price data, per the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 689.0357
$ws.Range("I80").Value = 846.8
$ws.Range("J80").Value = 507
$ws.Range("K80").Value = 2540.4
$ws.Range("L80").Value = 1521
$ws.Range("M80").Value = -1542.4
$ws.Range("N80").Value = -3517

$ws.Range("H83").Value = 689.0357
$ws.Range("I83").Value = 846.8
$ws.Range("J83").Value = 507
$ws.Range("K83").Value = 7621.2
$ws.Range("L83").Value = 4563
$ws.Range("M83").Value = -2629.2
$ws.Range("N83").Value = -14547

$ws.Range("H116").Value = 4299.737
$ws.Range("I116").Value = 2482.9167
$ws.Range("J116").Value = 7414.2856
$ws.Range("K116").Value = 2482.9167
$ws.Range("L116").Value = 7414.2856
$ws.Range("M116").Value = 959.0832999999998
$ws.Range("N116").Value = -14298.2856

$ws.Range("H138").Value = 1878
$ws.Range("I138").Value = 825.2195
$ws.Range("J138").Value = 5475
$ws.Range("K138").Value = 2475.6585
$ws.Range("L138").Value = 16425
$ws.Range("M138").Value = 2664.3415
$ws.Range("N138").Value = -26705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2115.04
$ws.Range("I32").Value = 1853.337
$ws.Range("J32").Value = 4232.4546
$ws.Range("K32").Value = 1853.337
$ws.Range("L32").Value = 4232.4546
$ws.Range("M32").Value = -1566.337
$ws.Range("N32").Value = -4806.4546

$ws.Range("H95").Value = 22334.834
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 22334.834
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 22334.834
$ws.Range("N95").Value = -27826.834

$ws.Range("H96").Value = 27448
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 27448
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 27448
$ws.Range("N96").Value = -32940

$ws.Range("H110").Value = 1696.7
$ws.Range("I110").Value = 1427.125
$ws.Range("J110").Value = 2775
$ws.Range("K110").Value = 1427.125
$ws.Range("L110").Value = 2775
$ws.Range("M110").Value = 617.875
$ws.Range("N110").Value = -6865

$ws.Range("H132").Value = 4334.064
$ws.Range("I132").Value = 2003.0476
$ws.Range("J132").Value = 6216.8076
$ws.Range("K132").Value = 6009.142800000001
$ws.Range("L132").Value = 18650.4228
$ws.Range("M132").Value = -3479.142800000001
$ws.Range("N132").Value = -23710.4228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 10000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = $null
$ws.Range("N36").Value = -10776

$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -10320

$ws.Range("H86").Value = 43482816
$ws.Range("I86").Value = 66670500
$ws.Range("J86").Value = 5901
$ws.Range("K86").Value = 66670500
$ws.Range("L86").Value = 5901
$ws.Range("M86").Value = -66669377
$ws.Range("N86").Value = -8147

$ws.Range("H89").Value = 43482816
$ws.Range("I89").Value = 66670500
$ws.Range("J89").Value = 5901
$ws.Range("K89").Value = 333352500
$ws.Range("L89").Value = 29505
$ws.Range("M89").Value = -333346884
$ws.Range("N89").Value = -40737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 197.36363
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 358.83334
$ws.Range("K12").Value = 10.8
$ws.Range("L12").Value = 1076.50002
$ws.Range("M12").Value = 162.2
$ws.Range("N12").Value = -1422.50002

$ws.Range("H21").Value = 1750
$ws.Range("I21").Value = 650
$ws.Range("J21").Value = 2850
$ws.Range("K21").Value = 1950
$ws.Range("L21").Value = 8550
$ws.Range("M21").Value = -1777
$ws.Range("N21").Value = -8896

$ws.Range("H25").Value = 83337840
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 83337840
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 250013520
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = -250013858

$ws.Range("H30").Value = 83337840
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 83337840
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 250013520
$ws.Range("M30").Value = $null
$ws.Range("N30").Value = -250013724

$ws.Range("H47").Value = 456.5
$ws.Range("I47").Value = 138.4
$ws.Range("J47").Value = 986.6667
$ws.Range("K47").Value = 415.2
$ws.Range("L47").Value = 2960.0001
$ws.Range("M47").Value = 15.79999999999995
$ws.Range("N47").Value = -3822.0001

$ws.Range("H131").Value = 371178.84
$ws.Range("I131").Value = 769669.7
$ws.Range("J131").Value = 1151.6428
$ws.Range("K131").Value = 2309009.1
$ws.Range("L131").Value = 3454.9284
$ws.Range("M131").Value = -2303969.1
$ws.Range("N131").Value = -13534.9284

$ws.Range("H137").Value = 3790.8667
$ws.Range("I137").Value = 2821
$ws.Range("J137").Value = 4275.8
$ws.Range("K137").Value = 8463
$ws.Range("L137").Value = 12827.4
$ws.Range("M137").Value = -3363
$ws.Range("N137").Value = -23027.4

$ws.Range("H139").Value = 2209.9443
$ws.Range("I139").Value = 964.3889
$ws.Range("J139").Value = 3455.5
$ws.Range("K139").Value = 2893.1667
$ws.Range("L139").Value = 10366.5
$ws.Range("M139").Value = 2246.8333
$ws.Range("N139").Value = -20646.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2289.0715
$ws.Range("I80").Value = 2283.5
$ws.Range("J80").Value = 2303
$ws.Range("K80").Value = 2283.5
$ws.Range("L80").Value = 2303
$ws.Range("M80").Value = -1285.5
$ws.Range("N80").Value = -4299

$ws.Range("H83").Value = 2289.0715
$ws.Range("I83").Value = 2283.5
$ws.Range("J83").Value = 2303
$ws.Range("K83").Value = 11417.5
$ws.Range("L83").Value = 11515
$ws.Range("M83").Value = -6425.5
$ws.Range("N83").Value = -21499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 38502936
$ws.Range("I132").Value = 66735788
$ws.Range("J132").Value = 3588
$ws.Range("K132").Value = 200207364
$ws.Range("L132").Value = 10764
$ws.Range("M132").Value = -200204834
$ws.Range("N132").Value = -15824

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3149.818
$ws.Range("I122").Value = 2558.5278
$ws.Range("J122").Value = 5810.625
$ws.Range("K122").Value = 7675.5834
$ws.Range("L122").Value = 17431.875
$ws.Range("M122").Value = -5225.5834
$ws.Range("N122").Value = -22331.875

$ws.Range("H132").Value = 3044
$ws.Range("I132").Value = 3239
$ws.Range("J132").Value = 2929.2942
$ws.Range("K132").Value = 9717
$ws.Range("L132").Value = 8787.882599999999
$ws.Range("M132").Value = -7187
$ws.Range("N132").Value = -13847.8826
